# Applies the cryptos list price/volume refresh described in the commit
# "Updated cryptos list on Fri May 24 17:36:17 UTC 2024 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.669.58"
$ws.Range("E2").Value = "  +0.86%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.700.56"
$ws.Range("E3").Value = "  -3.19%  "

# Row 4
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.48"
$ws.Range("E5").Value = "  +0.31%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.15"
$ws.Range("E6").Value = "  -4.26%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.700.14"
$ws.Range("E7").Value = "  -3.13%  "

# Row 8
$ws.Range("E8").Value = "  -0.09%  "

# Row 9
$ws.Range("E9").Value = "  +0.98%  "

# Row 10
$ws.Range("E10").Value = "  +2.52%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.28"
$ws.Range("E11").Value = "  -0.29%  "

# Row 12
$ws.Range("E12").Value = "  -1.52%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.04"
$ws.Range("E13").Value = "  -0.11%  "

# Row 14
$ws.Range("E14").Value = "  -0.87%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.313.82"
$ws.Range("E15").Value = "  -3.09%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.693.92"
$ws.Range("E16").Value = "  -2.96%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.629.64"
$ws.Range("E17").Value = "  +0.71%  "

# Row 18
$ws.Range("E18").Value = "  +1.19%  "

# Row 19
$ws.Range("E19").Value = "  -0.66%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.15"
$ws.Range("E20").Value = "  +5.37%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "492.95"
$ws.Range("E21").Value = "  +0.49%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.13"
$ws.Range("E22").Value = "  -1.13%  "

# Row 23
$ws.Range("E23").Value = "  -1.84%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.35"
$ws.Range("E24").Value = "  -0.62%  "

# Row 25
$ws.Range("B25").Value = "Fetch.AI"
$ws.Range("C25").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.29"
$ws.Range("E25").Value = "  -4.16%  "

# Row 26
$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000141"
$ws.Range("E26").Value = "  +1.53%  "

# Row 27
$ws.Range("E27").Value = "  -1.13%  "

# Row 28
$ws.Range("E28").Value = "  -2.01%  "

# Row 29
$ws.Range("E29").Value = "  +0.05%  "

# Row 30
$ws.Range("E30").Value = "  -0.45%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.81"
$ws.Range("E31").Value = "  +0.70%  "

# Row 32
$ws.Range("E32").Value = "  -2.77%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.44"
$ws.Range("E33").Value = "  -4.19%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.831.61"
$ws.Range("E34").Value = "  -3.19%  "

# Row 35
$ws.Range("E35").Value = "  -0.88%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.635.72"
$ws.Range("E36").Value = "  -3.22%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.19%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.998"
$ws.Range("E38").Value = "  -1.34%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.73"
$ws.Range("E39").Value = "  -1.00%  "

# Row 40
$ws.Range("E40").Value = "  -3.83%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.322"
$ws.Range("E41").Value = "  -1.47%  "

# Row 42
$ws.Range("E42").Value = "  +0.10%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "431.57"
$ws.Range("E43").Value = "  -4.45%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.97"
$ws.Range("E44").Value = "  -1.21%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.82"
$ws.Range("E45").Value = "  -2.35%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.36"
$ws.Range("E46").Value = "  +0.88%  "

# Row 47
$ws.Range("E47").Value = "  +0.01%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.20"
$ws.Range("E48").Value = "  -3.13%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.66"
$ws.Range("E49").Value = "  +2.06%  "

# Row 50
$ws.Range("E50").Value = "  -1.04%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.731.30"
